$d = $word.ActiveDocument

# --- Change 1: resize the first "Bihevioralni pogled_Trener" image ---
# The picture is a legacy VML <w:pict> shape, so it is not exposed via the
# InlineShapes/Shapes collections; edit its markup directly by replacing the
# exact paragraph Range that contains it with the same paragraph/run
# formatting, only changing the shape's height from 239.6pt to 239.05pt.
$imgParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $x = $p.Range.WordOpenXML
    if ($x.IndexOf("_x0000_i1026") -ge 0) {
        $imgParaIndex = $i
        break
    }
}

if ($imgParaIndex -gt 0) {
    $imgPara = $d.Paragraphs.Item($imgParaIndex)
    $imgRange = $imgPara.Range
    $newParaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:o="urn:schemas-microsoft-com:office:office"><w:body><w:p w:rsidR="00D426D6" w:rsidRDefault="00471A98" w:rsidP="001607C3"><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:pict><v:shape id="_x0000_i1026" type="#_x0000_t75" style="width:466pt;height:239.05pt"><v:imagedata r:id="rId6" o:title="Bihevioralni pogled_Trener"/></v:shape></w:pict></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    [void]$imgRange.InsertXML($newParaXml)
}

# --- Change 2: remove the "NodeJS MVC" paragraph entirely ---
# It is the whole paragraph that follows the
# "Specifikacija biblioteka i programskih okvira" heading and precedes the
# "RabbitMQ" paragraph; delete the whole paragraph (including its mark) so
# the document collapses back to the original flow.
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.IndexOf("NodeJS") -ge 0) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    $target = $d.Paragraphs.Item($targetIndex)
    [void]$target.Range.Delete()
}
